$d = $word.ActiveDocument

# Locate the paragraph containing the "Rpta" answer about array default values.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Se crearan 5 posiciones*") {
        $target = $i
        break
    }
}
if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$p = $d.Paragraphs.Item($target)
$paraStart = $p.Range.Start
$paraEnd = $p.Range.End

# Within this paragraph, find the lone "." that sits right after the
# (hidden) _GoBack bookmark.
$searchRange = $d.Range($paraStart, $paraEnd)
$found = $searchRange.Find.Execute(".", $false, $false, $false, $false, $false, $true, 0, $false, "", 0)
if (-not $found) {
    throw "Could not find trailing period"
}
$dotStart = $searchRange.Start
$dotEnd = $searchRange.End

# Merge it into the preceding text run (" ... valor "0"" -> " ... valor "0".") ...
$insertPoint = $d.Range($dotStart, $dotStart)
$insertPoint.InsertBefore(".")

# ... then remove the now-redundant standalone "." run (shifted by the
# character we just inserted).
$oldDot = $d.Range($dotStart + 1, $dotEnd + 1)
$oldDot.Delete()

# Collapse the block of blank / tab-only paragraphs that follow this
# paragraph, keeping only the last two of them.
$trailing = 0
$scan = $target + 1
while ($scan -le $d.Paragraphs.Count -and $d.Paragraphs.Item($scan).Range.Text.Trim() -eq "") {
    $trailing = $trailing + 1
    $scan = $scan + 1
}
$toDelete = $trailing - 2
if ($toDelete -gt 0) {
    $delStart = $d.Paragraphs.Item($target + 1).Range.Start
    $delEnd = $d.Paragraphs.Item($target + 1 + $toDelete).Range.Start
    $trim = $d.Range($delStart, $delEnd)
    $trim.Delete()
}

Write-Output ("Done. ParaCount=" + $d.Paragraphs.Count)
